# Auto-generated edit script: updates cryptos list values (Price / Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.461.98'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.628.20'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.93'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.646'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +5.37%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.78'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.389'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.71%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.41'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000186'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -6.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.101.73'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.280.51'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.622.93'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.23'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.67'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.36'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '344.95'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.75'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000113'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.33'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.56%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '556.51'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.162'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.90'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.34%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.43'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.26'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.410'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.99'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.65%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.92'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '154.03'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.43'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '158.38'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.97'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0596'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.68'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.634'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.102'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0250'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.07'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0240'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.42%  '
